# Atualização automática SALDO_PECAS (14/11/2025 19:26)
# Row 6: UF changed DF -> GO, FRU changed 1234567 -> 00P0098,
# DESCRICAO/MAQUINAS changed TESTE/DS8K -> A/A, CLIENTE/DATA_FIM/SLA updated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "GO"
$ws.Range("B6").Value = "00P0098"
$ws.Range("F6").Value = "A"
$ws.Range("G6").Value = "A"
$ws.Range("H6").Value = "A - (A 01/11/25_12H) - GO"

# Force text format on I6 before assigning so the dd/mm/yy-looking
# string is kept as plain text instead of being auto-parsed as a date.
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "01/11/25"

$ws.Range("J6").Value = "12H"
